# Minor fix on template:
#  - fix the repo URL in the "Software source:" shared string (B4 on the
#    "Statistics" sheet)
#  - turn that cell into a real hyperlink (adds the Hyperlink style/font +
#    the <hyperlinks> part + relationship)
#  - move the active selection from B2 to B5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistics")

# Correct the software-source URL text (shared string is reused in place).
$ws.Range("B4").Value = "https://github.com/auino/cryptocurrencies2excel"

# Turn B4 into a clickable hyperlink pointing at the corrected URL. This
# also creates/applies the builtin "Hyperlink" cell style (underline, theme
# color 10) the same way Excel does when you insert a hyperlink.
$ws.Hyperlinks.Add($ws.Range("B4"), $ws.Range("B4").Text) | Out-Null

# Update the sheet's active selection.
$ws.Range("B5").Select() | Out-Null
